$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete risk rows (old rows 11 and 12) that are no longer
# part of the refreshed secure test run.
$ws.Rows(11).Delete()
$ws.Rows(11).Delete()

# Refresh the remaining risk rows (2-10) with the updated risk catalogue:
# RAA % values recalculated and several risks replaced by the new set
# (Missing Build Infrastructure, Missing Cloud Hardening, Missing Hardening
# for app/database) while dropping CSRF, Identity Propagation and WAF rows.

# Row 2
$ws.Range("A2").Value2 = 'Elevated'
$ws.Range("B2").Value2 = 'Likely'
$ws.Range("C2").Value2 = 'Medium'
$ws.Range("D2").Value2 = 'Tampering'
$ws.Range("E2").Value2 = 'Development'
$ws.Range("F2").Value2 = 'CWE-79'
$ws.Range("G2").Value2 = 'Cross-Site Scripting (XSS)'
$ws.Range("H2").Value2 = 'app'
$ws.Range("I2").Value2 = ''
$ws.Range("J2").Value2 = 58
$ws.Range("K2").Value2 = 'Cross-Site Scripting (XSS) risk at app'
$ws.Range("L2").Value2 = 'XSS Prevention'
$ws.Range("M2").Value2 = 'Try to encode all values sent back to the browser and also handle DOM-manipulations in a safe way to avoid DOM-based XSS. When a third-party product is used instead of custom developed software, check if the product applies the proper mitigation and ensure a reasonable patch-level.'
$ws.Range("N2").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O2").Value2 = 'cross-site-scripting@app'
$ws.Range("P2").Value2 = 'Unchecked'

# Row 3
$ws.Range("A3").Value2 = 'Elevated'
$ws.Range("B3").Value2 = 'Likely'
$ws.Range("C3").Value2 = 'Medium'
$ws.Range("D3").Value2 = 'Tampering'
$ws.Range("E3").Value2 = 'Development'
$ws.Range("F3").Value2 = 'CWE-79'
$ws.Range("G3").Value2 = 'Cross-Site Scripting (XSS)'
$ws.Range("H3").Value2 = 'webapp'
$ws.Range("I3").Value2 = ''
$ws.Range("J3").Value2 = 12
$ws.Range("K3").Value2 = 'Cross-Site Scripting (XSS) risk at webapp'
$ws.Range("L3").Value2 = 'XSS Prevention'
$ws.Range("M3").Value2 = 'Try to encode all values sent back to the browser and also handle DOM-manipulations in a safe way to avoid DOM-based XSS. When a third-party product is used instead of custom developed software, check if the product applies the proper mitigation and ensure a reasonable patch-level.'
$ws.Range("N3").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O3").Value2 = 'cross-site-scripting@webapp'
$ws.Range("P3").Value2 = 'Unchecked'

# Row 4
$ws.Range("A4").Value2 = 'Elevated'
$ws.Range("B4").Value2 = 'Very Likely'
$ws.Range("C4").Value2 = 'Medium'
$ws.Range("D4").Value2 = 'Tampering'
$ws.Range("E4").Value2 = 'Development'
$ws.Range("F4").Value2 = 'CWE-89'
$ws.Range("G4").Value2 = 'SQL/NoSQL-Injection'
$ws.Range("H4").Value2 = 'app'
$ws.Range("I4").Value2 = 'to-database'
$ws.Range("J4").Value2 = 58
$ws.Range("K4").Value2 = 'SQL/NoSQL-Injection risk at app against database database via to-database'
$ws.Range("L4").Value2 = 'SQL/NoSQL-Injection Prevention'
$ws.Range("M4").Value2 = 'Try to use parameter binding to be safe from injection vulnerabilities. When a third-party product is used instead of custom developed software, check if the product applies the proper mitigation and ensure a reasonable patch-level.'
$ws.Range("N4").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O4").Value2 = 'sql-nosql-injection@app@database@app>to-database'
$ws.Range("P4").Value2 = 'Unchecked'

# Row 5
$ws.Range("A5").Value2 = 'Medium'
$ws.Range("B5").Value2 = 'Unlikely'
$ws.Range("C5").Value2 = 'Medium'
$ws.Range("D5").Value2 = 'Tampering'
$ws.Range("E5").Value2 = 'Architecture'
$ws.Range("F5").Value2 = 'CWE-1127'
$ws.Range("G5").Value2 = 'Missing Build Infrastructure'
$ws.Range("H5").Value2 = 'app'
$ws.Range("I5").Value2 = ''
$ws.Range("J5").Value2 = 58
$ws.Range("K5").Value2 = 'Missing Build Infrastructure in the threat model (referencing asset app as an example)'
$ws.Range("L5").Value2 = 'Build Pipeline Hardening'
$ws.Range("M5").Value2 = 'Include the build infrastructure in the model.'
$ws.Range("N5").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O5").Value2 = 'missing-build-infrastructure@app'
$ws.Range("P5").Value2 = 'Unchecked'

# Row 6
$ws.Range("A6").Value2 = 'Medium'
$ws.Range("B6").Value2 = 'Unlikely'
$ws.Range("C6").Value2 = 'High'
$ws.Range("D6").Value2 = 'Tampering'
$ws.Range("E6").Value2 = 'Operations'
$ws.Range("F6").Value2 = 'CWE-1008'
$ws.Range("G6").Value2 = 'Missing Cloud Hardening'
$ws.Range("H6").Value2 = ''
$ws.Range("I6").Value2 = ''
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = 'Missing Cloud Hardening risk at database'
$ws.Range("L6").Value2 = 'Cloud Hardening'
$ws.Range("M6").Value2 = 'Apply hardening of all cloud components and services, taking special care to follow the individual risk descriptions (which depend on the cloud provider tags in the model). <br><br>For <b>Amazon Web Services (AWS)</b>: Follow the <i>CIS Benchmark for Amazon Web Services</i> (see also the automated checks of cloud audit tools like <i>"PacBot", "CloudSploit", "CloudMapper", "ScoutSuite", or "Prowler AWS CIS Benchmark Tool"</i>). <br>For EC2 and other servers running Amazon Linux, follow the <i>CIS Benchmark for Amazon Linux</i> and switch to IMDSv2. <br>For S3 buckets follow the <i>Security Best Practices for Amazon S3</i> at <a href="https://docs.aws.amazon.com/AmazonS3/latest/dev/security-best-practices.html">https://docs.aws.amazon.com/AmazonS3/latest/dev/security-best-practices.html</a> to avoid accidental leakage. <br>Also take a look at some of these tools: <a href="https://github.com/toniblyx/my-arsenal-of-aws-security-tools">https://github.com/toniblyx/my-arsenal-of-aws-security-tools</a> <br><br>For <b>Microsoft Azure</b>: Follow the <i>CIS Benchmark for Microsoft Azure</i> (see also the automated checks of cloud audit tools like <i>"CloudSploit" or "ScoutSuite"</i>).<br><br>For <b>Google Cloud Platform</b>: Follow the <i>CIS Benchmark for Google Cloud Computing Platform</i> (see also the automated checks of cloud audit tools like <i>"CloudSploit" or "ScoutSuite"</i>). <br><br>For <b>Oracle Cloud Platform</b>: Follow the hardening best practices (see also the automated checks of cloud audit tools like <i>"CloudSploit"</i>).'
$ws.Range("N6").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O6").Value2 = 'missing-cloud-hardening@database'
$ws.Range("P6").Value2 = 'Unchecked'

# Row 7
$ws.Range("A7").Value2 = 'Medium'
$ws.Range("B7").Value2 = 'Likely'
$ws.Range("C7").Value2 = 'Low'
$ws.Range("D7").Value2 = 'Tampering'
$ws.Range("E7").Value2 = 'Operations'
$ws.Range("F7").Value2 = 'CWE-16'
$ws.Range("G7").Value2 = 'Missing Hardening'
$ws.Range("H7").Value2 = 'app'
$ws.Range("I7").Value2 = ''
$ws.Range("J7").Value2 = 58
$ws.Range("K7").Value2 = 'Missing Hardening risk at app'
$ws.Range("L7").Value2 = 'System Hardening'
$ws.Range("M7").Value2 = 'Try to apply all hardening best practices (like CIS benchmarks, OWASP recommendations, vendor recommendations, DevSec Hardening Framework, DBSAT for Oracle databases, and others).'
$ws.Range("N7").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O7").Value2 = 'missing-hardening@app'
$ws.Range("P7").Value2 = 'Unchecked'

# Row 8
$ws.Range("A8").Value2 = 'Medium'
$ws.Range("B8").Value2 = 'Likely'
$ws.Range("C8").Value2 = 'Low'
$ws.Range("D8").Value2 = 'Tampering'
$ws.Range("E8").Value2 = 'Operations'
$ws.Range("F8").Value2 = 'CWE-16'
$ws.Range("G8").Value2 = 'Missing Hardening'
$ws.Range("H8").Value2 = 'database'
$ws.Range("I8").Value2 = ''
$ws.Range("J8").Value2 = 100
$ws.Range("K8").Value2 = 'Missing Hardening risk at database'
$ws.Range("L8").Value2 = 'System Hardening'
$ws.Range("M8").Value2 = 'Try to apply all hardening best practices (like CIS benchmarks, OWASP recommendations, vendor recommendations, DevSec Hardening Framework, DBSAT for Oracle databases, and others).'
$ws.Range("N8").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O8").Value2 = 'missing-hardening@database'
$ws.Range("P8").Value2 = 'Unchecked'

# Row 9
$ws.Range("A9").Value2 = 'Medium'
$ws.Range("B9").Value2 = 'Unlikely'
$ws.Range("C9").Value2 = 'Medium'
$ws.Range("D9").Value2 = 'Spoofing'
$ws.Range("E9").Value2 = 'Architecture'
$ws.Range("F9").Value2 = 'CWE-287'
$ws.Range("G9").Value2 = 'Missing Identity Store'
$ws.Range("H9").Value2 = 'webapp'
$ws.Range("I9").Value2 = ''
$ws.Range("J9").Value2 = 12
$ws.Range("K9").Value2 = 'Missing Identity Store in the threat model (referencing asset webapp as an example)'
$ws.Range("L9").Value2 = 'Identity Store'
$ws.Range("M9").Value2 = 'Include an identity store in the model if the application has a login.'
$ws.Range("N9").Value2 = 'Are recommendations from the linked cheat sheet and referenced ASVS chapter applied?'
$ws.Range("O9").Value2 = 'missing-identity-store@webapp'
$ws.Range("P9").Value2 = 'Unchecked'

# Row 10
$ws.Range("A10").Value2 = 'Medium'
$ws.Range("B10").Value2 = 'Unlikely'
$ws.Range("C10").Value2 = 'Medium'
$ws.Range("D10").Value2 = 'Information Disclosure'
$ws.Range("E10").Value2 = 'Architecture'
$ws.Range("F10").Value2 = 'CWE-522'
$ws.Range("G10").Value2 = 'Missing Vault (Secret Storage)'
$ws.Range("H10").Value2 = 'app'
$ws.Range("I10").Value2 = ''
$ws.Range("J10").Value2 = 58
$ws.Range("K10").Value2 = 'Missing Vault (Secret Storage) in the threat model (referencing asset app as an example)'
$ws.Range("L10").Value2 = 'Vault (Secret Storage)'
$ws.Range("M10").Value2 = 'Consider using a Vault (Secret Storage) to securely store and access config secrets (like credentials, private keys, client certificates, etc.).'
$ws.Range("N10").Value2 = 'Is a Vault (Secret Storage) in place?'
$ws.Range("O10").Value2 = 'missing-vault@app'
$ws.Range("P10").Value2 = 'Unchecked'
